$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers: I1 = "I0", J1 = "IF" (reuse header style from H1) ---
$ws.Cells.Item(1, 8).Copy()
$ws.Cells.Item(1, 9).PasteSpecial(-4122)
$ws.Cells.Item(1, 10).PasteSpecial(-4122)
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# --- Data rows: column I (I0) and column J (IF) ---
$I0 = @{2=1;3=1;4=1;5=1;6=1;7=1;8=1;9=1;10=1;11=1;12=1;13=1;14=1;15=1;16=1;17=1;18=1;19=1;20=1;21=1;22=1;23=1;24=1;25=1;26=1;27=1;28=6;29=1;30=4;31=1}
$IF = @{2=8;3=6;4=6;5=5;6=6;7=6;8=5;9=5;10=6;11=5;12=3;13=6;14=4;15=6;16=4;17=4;18=6;19=5;20=5;21=6;22=4;23=5;24=5;25=6;26=5;27=6;28=9;29=4;30=6;31=2}

foreach ($r in 2..31) {
    $ws.Cells.Item($r, 9).Value = $I0[$r]
    $ws.Cells.Item($r, 10).Value = $IF[$r]
}
